$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.234.21"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.430.12"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Formula = "'407.25"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Formula = "'134.03"
$ws.Range("E6").Value = "  +3.77%  "
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Formula = "'0.688"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Formula = "'0.123"
$ws.Range("E10").Value = "  -5.36%  "
$ws.Range("D11").Formula = "'42.26"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Formula = "'8.45"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Formula = "'19.96"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "3.417.16"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "62.206.14"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Formula = "'11.32"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Formula = "'1.02"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").Formula = "'3.20"
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("D21").Formula = "'84.37"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("D22").Formula = "'314.07"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").Formula = "'12.92"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Formula = "'3.15"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Formula = "'4.77"
$ws.Range("E25").Value = "  +8.79%  "
$ws.Range("D26").Formula = "'29.72"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Formula = "'2.79"
$ws.Range("E28").Value = "  +5.86%  "
$ws.Range("D29").Formula = "'7.56"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Formula = "'0.174"
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("D32").Formula = "'42.48"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Formula = "'11.39"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("D35").Formula = "'0.0486"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").Formula = "'51.46"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Formula = "'0.312"
$ws.Range("E40").Value = "  +9.34%  "
$ws.Range("D41").Formula = "'138.54"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Formula = "'0.125"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Formula = "'1.99"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Formula = "'4.04"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Formula = "'16.80"
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("D46").Formula = "'2.23"
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("D47").Formula = "'21.32"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("D48").Value = "2.121.26"
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("D49").Formula = "'2.30"
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("D50").Formula = "'1.93"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("E51").Value = "  +17.97%  "
